$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove Kiko row; Ana and Hugo rows shift up by one.
$ws.Rows(4).Delete()

$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "589.0"
$ws.Range("E5").Value = "25/12/2020"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "89.0"
$ws.Range("G5").Value = "25/12/2020"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "89.0"
$ws.Range("I5").Value = "25/12/2020"
$ws.Range("J5").NumberFormat = "@"
$ws.Range("J5").Value = "89.0"
$ws.Range("K5").Value = "25/12/2020"
$ws.Range("L5").NumberFormat = "@"
$ws.Range("L5").Value = "89.0"
$ws.Range("M5").Value = "25/12/2020"
$ws.Range("N5").NumberFormat = "@"
$ws.Range("N5").Value = "89.0"
$ws.Range("O5").Value = "25/12/2020"
$ws.Range("P5").NumberFormat = "@"
$ws.Range("P5").Value = "89.0"
$ws.Range("Q5").Value = "25/12/2020"
$ws.Range("R5").NumberFormat = "@"
$ws.Range("R5").Value = "89.0"
$ws.Range("S5").Value = "25/12/2020"
$ws.Range("T5").NumberFormat = "@"
$ws.Range("T5").Value = "89.0"
$ws.Range("U5").Value = "25/12/2020"
$ws.Range("V5").NumberFormat = "@"
$ws.Range("V5").Value = "89.0"
$ws.Range("W5").Value = "25/12/2020"
$ws.Range("X5").NumberFormat = "@"
$ws.Range("X5").Value = "89.0"
$ws.Range("Y5").Value = "25/12/2020"
$ws.Range("Z5").NumberFormat = "@"
$ws.Range("Z5").Value = "89.0"
$ws.Range("A6").Value = 400
$ws.Range("B6").Value = "Moises"
$ws.Range("C6").Value = "24/12/2020"
$ws.Range("D6").Value = 400
$ws.Range("A7").Value = 350
$ws.Range("B7").Value = "Joaquim"
$ws.Range("C7").Value = "24/12/2020"
$ws.Range("D7").Value = 350
$ws.Range("A8").Value = 300
$ws.Range("B8").Value = "Kirliaa"
$ws.Range("C8").Value = "24/12/2020"
$ws.Range("D8").Value = 300
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "880.0"
$ws.Range("B9").Value = "Juka"
$ws.Range("C9").Value = "24/12/2020"
$ws.Range("D9").Value = 800
$ws.Range("E9").Value = "25/12/2020"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "80.0"
$ws.Range("G9").Value = "25/12/2020"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "80.0"
$ws.Range("I9").Value = "25/12/2020"
$ws.Range("J9").NumberFormat = "@"
$ws.Range("J9").Value = "80.0"
$ws.Range("K9").Value = "25/12/2020"
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "80.0"
$ws.Range("M9").Value = "25/12/2020"
$ws.Range("N9").NumberFormat = "@"
$ws.Range("N9").Value = "80.0"
$ws.Range("O9").Value = "25/12/2020"
$ws.Range("P9").NumberFormat = "@"
$ws.Range("P9").Value = "80.0"
$ws.Range("Q9").Value = "25/12/2020"
$ws.Range("R9").NumberFormat = "@"
$ws.Range("R9").Value = "80.0"
$ws.Range("S9").Value = "25/12/2020"
$ws.Range("T9").NumberFormat = "@"
$ws.Range("T9").Value = "80.0"
$ws.Range("U9").Value = "25/12/2020"
$ws.Range("V9").NumberFormat = "@"
$ws.Range("V9").Value = "80.0"
$ws.Range("W9").Value = "25/12/2020"
$ws.Range("X9").NumberFormat = "@"
$ws.Range("X9").Value = "80.0"
$ws.Range("Y9").Value = "25/12/2020"
$ws.Range("Z9").NumberFormat = "@"
$ws.Range("Z9").Value = "80.0"
$ws.Range("AA9").Value = "25/12/2020"
$ws.Range("AB9").NumberFormat = "@"
$ws.Range("AB9").Value = "80.0"
$ws.Range("A10").Value = 900
$ws.Range("B10").Value = "Oseias"
$ws.Range("C10").Value = "24/12/2020"
$ws.Range("D10").Value = 900
$ws.Range("A11").Value = 890
$ws.Range("B11").Value = "Haas"
$ws.Range("C11").Value = "24/12/2020"
$ws.Range("D11").Value = 890
$ws.Range("A12").Value = 900
$ws.Range("B12").Value = "Test2"
$ws.Range("C12").Value = "24/12/2020"
$ws.Range("D12").Value = 900
$ws.Range("A13").Value = 670
$ws.Range("B13").Value = "Janete"
$ws.Range("C13").Value = "24/12/2020"
$ws.Range("D13").Value = 670
$ws.Range("A14").Value = 700
$ws.Range("B14").Value = "Cleide"
$ws.Range("C14").Value = "24/12/2020"
$ws.Range("D14").Value = 700
$ws.Range("A15").Value = 670
$ws.Range("B15").Value = "June"
$ws.Range("C15").Value = "24/12/2020"
$ws.Range("D15").Value = 670
$ws.Range("A16").Value = 789.65
$ws.Range("B16").Value = "Jurema"
$ws.Range("C16").Value = "25/12/2020"
$ws.Range("D16").Value = 789.65
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "85.0"
$ws.Range("B17").Value = "Joarez"
$ws.Range("C17").Value = "25/12/2020"
$ws.Range("D17").Value = 976
$ws.Range("E17").Value = "25/12/2020"
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value = "-267.0"
$ws.Range("G17").Value = "25/12/2020"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "-267.0"
$ws.Range("I17").Value = "25/12/2020"
$ws.Range("J17").NumberFormat = "@"
$ws.Range("J17").Value = "-267.0"
$ws.Range("K17").Value = "25/12/2020"
$ws.Range("L17").NumberFormat = "@"
$ws.Range("L17").Value = "-267.0"
$ws.Range("M17").Value = "25/12/2020"
$ws.Range("N17").NumberFormat = "@"
$ws.Range("N17").Value = "-267.0"
$ws.Range("O17").Value = "25/12/2020"
$ws.Range("P17").NumberFormat = "@"
$ws.Range("P17").Value = "-267.0"
$ws.Range("Q17").Value = "25/12/2020"
$ws.Range("R17").NumberFormat = "@"
$ws.Range("R17").Value = "-267.0"
$ws.Range("S17").Value = "25/12/2020"
$ws.Range("T17").NumberFormat = "@"
$ws.Range("T17").Value = "-267.0"
$ws.Range("U17").Value = "25/12/2020"
$ws.Range("V17").NumberFormat = "@"
$ws.Range("V17").Value = "-267.0"
$ws.Range("W17").Value = "25/12/2020"
$ws.Range("X17").NumberFormat = "@"
$ws.Range("X17").Value = "-267.0"
$ws.Range("Y17").Value = "25/12/2020"
$ws.Range("Z17").NumberFormat = "@"
$ws.Range("Z17").Value = "-267.0"
$ws.Range("AA17").Value = "25/12/2020"
$ws.Range("AB17").NumberFormat = "@"
$ws.Range("AB17").Value = "-267.0"
$ws.Range("AC17").Value = "25/12/2020"
$ws.Range("AD17").NumberFormat = "@"
$ws.Range("AD17").Value = "-267.0"
$ws.Range("AE17").Value = "25/12/2020"
$ws.Range("AF17").NumberFormat = "@"
$ws.Range("AF17").Value = "50.0"
$ws.Range("AG17").Value = "25/12/2020"
$ws.Range("AH17").NumberFormat = "@"
$ws.Range("AH17").Value = "-560.0"
$ws.Range("AI17").Value = "25/12/2020"
$ws.Range("AJ17").NumberFormat = "@"
$ws.Range("AJ17").Value = "-20.0"
$ws.Range("AK17").Value = "25/12/2020"
$ws.Range("AL17").NumberFormat = "@"
$ws.Range("AL17").Value = "-10.0"
$ws.Range("AM17").Value = "25/12/2020"
$ws.Range("AN17").NumberFormat = "@"
$ws.Range("AN17").Value = "-5.0"
